$wb = $excel.ActiveWorkbook

# --- Numeric single-cell updates (Week 13 stat increments) ---
$ws = $wb.Worksheets.Item("OFF")
$ws.Range("C2").Value = 284
$ws.Range("E2").Value = 16
$ws.Range("F2").Value = 124
$ws.Range("G2").Value = 113
$ws.Range("I2").Value = 18
$ws.Range("J2").Value = 66
$ws.Range("L2").Value = 568
$ws.Range("M2").Value = 375
$ws.Range("O2").Value = 24
$ws.Range("Q2").Value = 976
$ws.Range("C3").Value = 469
$ws.Range("E3").Value = 62
$ws.Range("F3").Value = 189
$ws.Range("G3").Value = 72
$ws.Range("H3").Value = 51
$ws.Range("I3").Value = 88
$ws.Range("J3").Value = 103
$ws.Range("N3").Value = 31

$ws = $wb.Worksheets.Item("DEF")
$ws.Range("C2").Value = 369
$ws.Range("E2").Value = 19
$ws.Range("F2").Value = 101
$ws.Range("G2").Value = 116
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 14
$ws.Range("J2").Value = 65
$ws.Range("L2").Value = 526
$ws.Range("M2").Value = 318
$ws.Range("O2").Value = 45
$ws.Range("P2").Value = 23
$ws.Range("Q2").Value = 988
$ws.Range("C3").Value = 342
$ws.Range("E3").Value = 77
$ws.Range("H3").Value = 55

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B2").Value = 183
$ws.Range("D2").Value = 83
$ws.Range("F2").Value = 94
$ws.Range("G2").Value = 92
$ws.Range("J2").Value = 36
$ws.Range("K2").Value = 32
$ws.Range("B3").Value = 117

$ws = $wb.Worksheets.Item("TURNS")
$ws.Range("D2").Value = 18
$ws.Range("E2").Value = 23

$ws = $wb.Worksheets.Item("PEN")
$ws.Range("B2").Value = 36
$ws.Range("B3").Value = 36

# --- Appended space-separated numbers to long play-log strings ---
$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B2").Value = '-3 12 6 9 2 4 16 0 -2 2 -4 -1 4 8 7 -1 3 3 0 12 -2 6 0 1 3 1 2 14 3 1 16 0 1 3 14 6 4 18 3 4 2 4 3 0 2 10 2 3 1 4 16 2 13 13 0 2 6 4 1 9 8 4 1 2 -3 6 0 5 17 7 -3 3 10 -3 11 3 6 4 4 4 4 0 3 -3 1 2 2 -5 3 2 -2 11 3 3 3 -3 -2 7 6 3 3 4 1 2 4 2 5 7 2 2 6 4 34 1 9 10 3 9 4 2 3 5 4 1 13 3 -3 3 6 -1 7 -6 2 0 6 3 7 14 17 12 5 -2 10 3 2 11 6 2 7 6 0 2 1 26 5 1 1 6 6 2 0 3 6 3 0 19 21 8 -2 2 4 6 1 2 7 6 6 18 4 2 2 5 9 17 2 12 3 -2 2 2 3 8 14 1 -1 6 7 7 -1 5 2 -6 1 3 2 0 1 5 3 1 1 2 2 2 1 -1 -3 9 13 15 0 2 8 1 14 0 -3 -1 12 6 3 6 -4 4 0 2 21 2 3 6 11 8 6 7 0 -1 8 3 1 11 8 0 8 -2 31 24 -4 5 12 12 4 5 7 3 0 -2 4 0 2 5 4 0 8 0 2 6 1 8 15 6 -1 2 2 6 -2 6 5 3 0 3 11 6 3 14 1 3 3 3 1 4 5 -4 4 13 6 3 2 -2 11 2 4 8 5 2 6 5 24 4 1 0 6 4 3 9 26 6 1 1 10 2 1 3 51 3 2 4 -1 5 3 11 5 2 3 11 4 4 1 3 3 2 2 11 6 2 5 4 2 3 12 5 -1 4 2 3 3 2 4 18 0 -3 1 0 7 18 3 2 1 1 7 5 3 12 6 -2 7 8 1 0 5 16 5 2 3 5 -5 9 0 2 6 0 13 4 3 9 -1 12 6 3 -2 -1 0 3 1 1 3 7 17 2 15 3 7 10 4 4 5 1 18 8 11 8 2 6 9 1 8 -3 1 3 1 3 3 3 7 3 2 2 11 2 6 15 25 2 3 46 4 6 1 10 1 -1 2 5 4 0 8 4 4 4 1 4 15 2 5 0 0 7 3 4 2 8 1 0 4 6 8 1 -2 4 4 8 2 3 4 3 2 2 0 3 0 7 6 5 1 6 9 1 2 2 3 4 3 1 8 0 17 8 -5 7 9 6 3 0 8 5 41 12 7 7 1 5 1 4 3 8 2 5 -1 4 0 8 22 2 11 9 9 8 2 -1 3 5 -1 0 0 3 1 4 12 4 6 6 0 8 3 1 3 14 3 7 3 -3 7 5 2 1 7 5 2 3 4 0 4 5 0 4 5 4 -1 6 34 1 0 -1 -2 6 5 7 13 2 6 2 12 3 -1 4 0 7 15 12 -3 1 22 4 -1 1 7 5 4 2 3 -1 15 0 -4 4 0 4 11 7 13 24 8 15 0 1 7 6 3 9 2 3 10 0 3 8 9 4 28 3 8 5 2 6 6 9 5 2 8 -1 5 3 3 -3 1 2 4 2 3 -1 4 10 1 11 15 0 -5 9 6 3 2 -2 3 0 4 0 2 2 1 4 4 5 3 3 1 1 5 3 6 3 2 21 3 0 17 5 1'

$ws = $wb.Worksheets.Item("YDS")
$ws.Range("B3").Value = '3 5 4 9 17 11 5 20 4 11 5 4 8 11 8 17 7 7 0 5 8 13 -1 14 10 11 11 4 22 2 15 5 25 15 26 8 5 1 36 4 11 46 22 16 15 12 9 8 6 47 24 6 21 1 46 4 12 5 39 9 5 21 29 3 34 12 13 0 4 7 23 5 22 6 19 17 3 2 21 15 6 26 7 10 18 6 11 7 9 1 5 8 12 5 7 12 14 14 49 9 14 9 19 15 16 9 3 15 14 12 12 8 5 9 6 10 7 -7 14 6 11 10 6 8 14 10 22 9 2 11 4 11 6 4 2 15 13 4 7 5 14 22 8 6 12 23 13 15 3 9 5 12 8 9 14 8 5 12 8 3 17 6 6 6 11 8 9 38 5 6 8 18 4 6 12 20 12 8 12 5 41 11 21 6 6 14 25 15 13 9 7 22 5 1 18 9 10 5 12 3 18 21 21 11 4 4 23 6 39 -2 20 25 11 33 7 6 10 31 12 7 5 11 9 14 2 5 9 5 6 13 11 0 22 5 21 7 7 -5 10 12 21 6 4 5 7 -3 17 4 21 5 2 3 8 14 2 20 15 14 7 1 4 4 6 8 11 6 44 2 9 3 31 5 5 6 12 35 22 5 8 9 8 8 12 4 7 18 20 18 9 12 4 6 11 23 3 10 12 3 28 11 2 13 8 -1 4 7 6 7 23 7 18 5 6 8 12 19 22 19 13 8 4 9 12 9 16 7 2 15 8 5 3 9 12 4 17 15 11 11 27 8 22 22 10 7 6 7 55 6 3 20 4 16 27 6 17 10 6 2 0 2 4 24 6 13 4 17 50 15 11 7 10 10 18 3 9 9 20 27 8 10 1 7 7 5 13 5 6 8 18 7 8 10 16 19 6 14 16 27 32 15 19 3 51 20 56 4 36 16 5 7 3 2 37 7 19 2 10 13 16 4 10 4 7 35 3 16 16 4 13 21 14 13 4 16 1 11 21 2 11 10 2 4 19 6 14 20 9 3 4 7 14 10 1 4 17 6 7 3 2 6 14 12 20 7 2 6 17 23 6 8 6 10 6 13 9 11 7 15 6 6 34 7 -2 6 5 10 16 13 4 7 12 10 5 37 10 8 3 7 9 9 12 7 4 6 12 -1 13 10 13 5 1 35 5 6 1 7 7 22 13 6 7 8 6 9 2 41 3 23 11 13 1 28 3 7 9 7 17 4 7 14 18 2 41 7 13 22 11 5 11 9 5 6 13 10 11 6 5 19 7 5 25 31 9 20 1 16 7 8 5 16 21 2 37 9 5 5 1 19 23 8 24 35 61 15 53 41 1 -2 11 5 16 16 8 9 12 12 6 4 3 11 6 7 3 3 21 19 2 14 5 17 12 29 1 5 9 11 18 9 15 3 4 1 6 6 31 12 7 20 1 11 8 1 20 3 15 7 8 6 15 4 10 6 3 3 8 8 17 11 14 19 13 1 11 -1 20 8 4 3 7 1 16 0 4 9 9 14 1 7 2 9 11 13 7 9 5 22 8 12 15 3 5 -2 28 1 7 19 5 14 3 11 17 15 4 7 7 12 3 6 8 6 6 16 57 12 23 31 49 25 8 36 20 43 7 11 9 7 20 2 18 3 6 14 12 7 6 10 6 7 16 7 31 1 9 14 9 7 7 26 1 6 12 7 9 4 17 4 7 14 5 21 4 5 6 28 23 -3 5 8 14 16 6 11 19 7 11 7 26 5 7 6'

$ws = $wb.Worksheets.Item("YDS")
$ws.Range("C2").Value = '3 5 -4 5 7 -2 6 5 2 3 8 4 0 2 6 2 1 3 1 8 9 3 14 9 2 1 0 11 0 0 4 5 3 9 2 5 2 7 7 4 0 7 0 12 0 4 14 14 6 3 2 2 -2 4 9 0 12 14 7 6 1 3 5 15 4 6 4 2 1 1 2 8 6 -3 -1 10 7 1 -2 1 -2 11 16 10 3 2 5 1 1 3 6 0 2 3 0 2 1 1 9 4 5 3 1 12 4 6 10 5 4 2 0 0 6 2 23 6 -3 7 5 3 9 3 4 2 -1 9 1 9 4 2 31 0 5 5 5 12 14 5 11 5 3 4 6 8 3 17 0 8 11 5 1 8 4 5 5 7 3 9 0 -3 2 -4 3 0 9 -4 2 8 11 4 7 2 7 0 8 2 -3 20 -1 1 5 9 2 2 13 -1 2 3 7 4 9 3 -2 2 2 5 2 2 4 2 6 -2 2 3 6 8 5 10 5 22 13 18 7 1 2 6 2 19 3 6 1 5 9 8 10 0 4 3 -3 0 8 4 -2 -3 4 5 -2 11 6 3 25 2 7 12 2 -7 7 5 1 13 4 1 12 -4 13 5 2 3 6 7 1 -1 4 28 4 15 20 8 7 3 5 4 -6 5 3 15 6 -2 0 1 2 -1 -1 33 -2 3 0 13 -2 3 1 -3 1 1 0 5 -2 16 8 4 4 8 1 0 0 1 3 2 12 3 18 -9 4 10 -1 3 3 3 3 1 7 2 5 5 3 1 0 5 1 0 6 2 10 4 13 1 2 3 0 1 14 2 10 2 4 4 0 3 3 14 5 6 7 9 8 8 0 3 2 3 29 -3 3 2 8 6 1 21 28 3 9 2 7 3 7 1 8 -2 2 2 3 2 0 3 5 1 -1 4 0 3 2 -7 5 1 7 16 1 3 2 2 16 6 3 2 4 12 -1 4 2 1 11 5 10 8 -2 0 1 -2 1 -3 3 3 2 4 3 29 20 3 33 5 12 10 9 -1 5 7 1 2 5 0 3 3 0 1 2 13 -4 9 5 2 15 1 1 6 -4 0 19 3 9 2 11 -2 5 -2 50 6 6 2 4 1 -3 3 9 0 5 11 1 2 5 2 -3 13 2 4 2 2 -4 2 0 9 2 1 2 25 -2 5 18 3 8 -1 6 18 -1 2 -1 6 2 1 3 2 8 3 4 -2 6 1 12 1 8 2 4 2 7 -1 4 1 -1 11 4 -4 7 3 1 7 4 7 6 0 1 0 5 1 10 3 4 0 12 0 -6 2 -1 2 4 1 3 2 4 2 5 4 2 10 5 2 3 15 3 6 4 1 9 5 0 -1 23 9 -1 -2 9 5 9 0 -1 2 4 3 76 2 -3 7 4 6 3 19 0 4 3 1 2 3 0 -1 13 1 7 5 5 5 -3 2 2 5 3 -2 12 3 1 1 0 2 6 2 0 3 1 7 6 3 2 3 -4 -1 15 1 -2 3 0 4 11 7 -1 5 1 6 2 8 2 0 -2 3 5 5 15 6 5 4 -4 5 -2 -2 5 8 5 13 0 0 0 5 1 -1 3 8 10 6 10 3 3 11 2 4 5 1 6 0 18 4 3 7 -1 2 2 3 5 40 10 9 0 3 2 1 19 6 5 18 9 3 -3 1 5 2 -2 3 4 8 2 0 0 2 2 1 4 -5 4 2 4 1 2 4 -1 -1 6 64 10 6 0 16 3 3 2 -2 -4 10 17 2 1 5 4 7 3 4 1 22 11 6 -1 4 1 10 5 5 -3 -5 3 -1 2 2'

$ws = $wb.Worksheets.Item("YDS")
$ws.Range("C3").Value = '0 9 13 30 4 4 21 5 -2 2 10 69 -5 2 9 7 6 7 1 8 18 17 2 13 13 26 2 4 14 2 -1 7 27 15 15 15 12 20 14 23 9 0 12 11 3 10 8 8 19 3 13 2 12 8 12 20 13 3 16 7 5 8 10 6 9 20 18 31 25 17 -6 17 8 31 16 23 9 -1 16 12 9 8 14 12 22 0 11 7 15 5 11 3 10 17 6 7 14 12 9 4 4 5 9 37 15 6 8 4 16 11 20 16 14 7 6 9 12 6 3 3 12 15 6 6 6 4 4 12 7 8 5 6 22 11 5 15 4 12 15 16 8 9 6 -3 21 -2 20 2 37 8 19 7 5 8 9 0 22 10 16 11 9 4 7 16 8 28 10 15 9 9 13 5 22 9 8 8 7 5 8 14 12 1 8 3 41 11 16 19 7 6 27 13 10 5 11 26 9 55 8 50 3 11 2 2 2 9 7 4 6 28 9 5 4 9 11 1 9 35 6 10 10 11 7 14 9 9 43 0 11 9 14 5 5 2 12 7 11 4 11 10 9 7 13 10 9 13 15 9 4 5 3 8 14 15 3 17 11 2 13 13 1 17 20 7 4 12 49 6 22 15 9 18 2 3 10 12 23 5 11 9 5 5 12 4 4 13 8 20 19 4 13 1 6 5 6 17 8 9 17 7 3 6 10 16 11 -2 8 5 5 9 6 6 19 6 -6 9 -2 6 4 5 10 7 9 7 3 7 8 22 6 6 10 11 6 8 7 22 3 3 1 10 8 5 32 5 17 13 25 5 6 2 16 11 31 7 12 13 23 7 6 5 -1 7 10 7 9 18 25 9 6 23 22 0 6 32 2 8 31 10 16 2 7 12 6 1 19 11 6 13 27 11 4 1 1 11 31 11 21 8 30 5 5 7 14 1 8 9 5 2 4 9 11 29 1 9 1 11 11 16 11 8 3 3 5 15 9 11 33 17 16 5 6 8 9 7 0 71 1 12 12 5 14 6 6 12 9 10 6 19 22 5 9 24 4 -2 5 24 14 1 13 7 6 9 5 18 12 4 5 14 4 11 16 11 3 5 8 0 4 2 7 8 4 5 6 0 11 17 10 73 14 7 4 15 1 11 15 37 3 8 7 8 5 8 12 3 6 12 12 6 15 11 9 7 8 5 8 10 6 26 8 12 8 11 5 4 3 11 8 17 6 8 5 5 4 1 8 6 4 4 11 11 5 48 7 8 8 8 4 13 20 6 6 14 22 5 8 9 23 7 8 0 14 15 22 12 8 3 8 5 1 -1 6 12 10 4 10 -2 11 4 6 12 -1 6 7 5 5 4 5 28 9 4 16 9 4 6 7 3 3 2 9 10 4 21 18 14 28 19 14 12 11 11 6 9 -7 21 18 8 -4 18 18 5 23 9 18 0 3 8 7 8 14 13 9 7 10 17 9 24 11 -3 7 6 3 7 12 7'

$ws = $wb.Worksheets.Item("ST")
$ws.Range("B6").Value = '31 -1 39 34 20 25 24 25 18 18 30 29 22 58 20 28 60 20 40 34 15 13 20 38 28 25 25 27 35 31 29 20 53 49 27 5 0 15 29 25 26 75 23 19 26 0 20 0 30 18 19 18 45 33 28 23 22 12 12 27 20 17 17 27'

$ws = $wb.Worksheets.Item("ST")
$ws.Range("D3").Value = '53 39 39 41 45 44 72 46 42 51 57 71 59 44 52 58 36 46 40 53 54 60 12 61 53 65 46 57 68 52 37 54 46 39 62 51 57 54 34 67 54 47 58 29 42 53 23 29 47 45 55 46 36 42 41 54 37 60 44 47 47 48 45 32 35 53 40 55 19 41 48 59 41 48 41 37 55 48 43 49 33 34 59'

$ws = $wb.Worksheets.Item("ST")
$ws.Range("D4").Value = '0 0 1 0 10 0 0 0 0 24 40 8 2 0 0 0 0 0 15 6 0 0 0 23 -8 0 0 0 16 0 0 0 0 0 0 5 0 0 0 12 0 0 14 0 0 0 0 0 0 0 0 5 0 4 0 0 0 0 6 30 30 0 0 0 0 2 0 8 0 -4 8 14 0 0 8 0 7 12 0 0 0 0 0'

$ws = $wb.Worksheets.Item("ST")
$ws.Range("D5").Value = '16 10 2 23 19 0 0 7 0 11 0 38 -1 0 0 24 0 19 7 0 15 0 0 0 0 11 0 10 0 18 0 4 0 0 4 22 0 0 0 0 0 0 0 15 0 0 7 7 0 7 6 0 12 0 7 0 0 0 0 84 0 0 0 0 0 8 0 0 0 1 0 0 0 20 0 0 0 13 6 33 0 0 5 0 9 6 12 -1 7 0 6 14 16 0 9 0 0 0 0 0 0 0 26 0 0 11 10 1 18 8 0 0 0 0 0'
